# Generate Report for Handoff
#
# Refreshes the localization-status report after a new handoff xliff
# generation run: the six "Ready for handoff" rows (source rows 7,8,9 and
# 11,12,13 on each per-language sheet) now carry their handoff priority
# ("ht") and the handoff/generate timestamps bump forward a few seconds.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 13)

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
# rows that just got a fresh handoff xliff.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-02 08:24:44"
}

# zh-cn sheet: Priority (E) gets the "ht" handoff-type tag, and the
# Latest Handoff Datetime (H) moves to the new generation time.
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-02 08:24:39"
}

# de-de sheet: same Priority tag, and its Latest Handoff Datetime lines
# up with the Overview's refreshed generate date.
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-02 08:24:44"
}
